# Fix field validation and UI guide in Excel upload
# Update test Excel file to match the standard template column layout.
#
# Old "Input" sheet columns (A:Q):
#   A:거래처명 B:현장명 C:발주일 D:납기일 E:발주번호 F:품목 G:규격 H:수량
#   I:단위 J:단가 K:공급가액 L:부가세 M:합계 N:대분류 O:중분류 P:소분류 Q:비고
#
# New "Input" sheet columns (A:P), standard template format:
#   A:발주일자 B:납기일자 C:거래처명 D:거래처 이메일 E:납품처명 F:납품처 이메일
#   G:프로젝트명 H:대분류 I:중분류 J:소분류 K:품목명 L:규격 M:수량 N:단가
#   O:총금액 P:비고

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the old trailing 비고 (Q) column; the new layout is only A:P wide and
# the replacement 비고 column lives at P (former 소분류 slot), which we
# overwrite explicitly below.
$ws.Columns.Item(17).Delete()

# --- Header row -------------------------------------------------------
$headers = New-Object 'object[,]' 1,16
$headers[0,0]  = "발주일자"
$headers[0,1]  = "납기일자"
$headers[0,2]  = "거래처명"
$headers[0,3]  = "거래처 이메일"
$headers[0,4]  = "납품처명"
$headers[0,5]  = "납품처 이메일"
$headers[0,6]  = "프로젝트명"
$headers[0,7]  = "대분류"
$headers[0,8]  = "중분류"
$headers[0,9]  = "소분류"
$headers[0,10] = "품목명"
$headers[0,11] = "규격"
$headers[0,12] = "수량"
$headers[0,13] = "단가"
$headers[0,14] = "총금액"
$headers[0,15] = "비고"
$ws.Range("A1:P1").Value = $headers
# Headers no longer carry the bold/border header style in the new template.
$ws.Range("A1:P1").Style = "Normal"

# --- Data rows ----------------------------------------------------------
# Columns: 발주일자 납기일자 거래처명 거래처이메일 납품처명 납품처이메일
#          프로젝트명 대분류 중분류 소분류 품목명 규격 수량 단가 총금액
# 발주일자/납기일자 are kept as plain text (leading apostrophe keeps Excel
# from auto-converting the "yyyy-mm-dd" strings into real date serials).
$rows = @(
    ,@("'2025-09-01","'2025-09-08","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","3) 기타","기타","스텐망 1168*343","KS규격-1",8,29000,255200)
    ,@("'2025-09-12","'2025-10-13","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","3) 기타","기타","스텐망 1023*1100","KS규격-2",2,29000,63800)
    ,@("'2025-09-17","'2025-09-29","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","3) 기타","기타","스텐망 1010*1900","KS규격-3",3,29000,95700)
    ,@("'2025-08-20","'2025-08-31","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","2월 운반비","KS규격-4",1,0,0)
    ,@("'2025-08-22","'2025-09-15","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","3) 기타","기타","스텐망 1088*1088","KS규격-5",2,29000,63800)
    ,@("'2025-09-15","'2025-10-16","케이에스파워텍","케이에스파워텍@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","2. 부자재비","3) 기타","기타","스텐망 1083*2145","KS규격-6",3,29000,95700)
)

$data = New-Object 'object[,]' 6,15
for ($r = 0; $r -lt 6; $r++) {
    for ($c = 0; $c -lt 15; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}
$ws.Range("A2:O7").Value = $data
# Strip the auto-applied "Text" number format from the date columns so the
# cells carry no style attribute, matching the template's plain cells.
$ws.Range("A2:B7").Style = "Normal"

# 비고 (P) only has content for the 2nd and 3rd shipment rows; all other
# rows must end up with no P cell at all.
$ws.Range("P2").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("P4").Value = "2차"
$ws.Range("P5").ClearContents()
$ws.Range("P6").ClearContents()
$ws.Range("P7").Value = "3차"

# --- 갑지 / 을지 sheets: drop the stray empty 비고 (I) cells ------------
# Rows 2,3,5,6 have no remark, so their I cell must be absent entirely
# (rows 4 and 7 keep their "2차"/"3차" remark).
for ($s = 2; $s -le 3; $s++) {
    $sheet = $wb.Worksheets.Item($s)
    $sheet.Range("I2").ClearContents()
    $sheet.Range("I3").ClearContents()
    $sheet.Range("I5").ClearContents()
    $sheet.Range("I6").ClearContents()
}

Write-Output "Input/갑지/을지 sheets updated to standard template layout"
